$d = $word.ActiveDocument

$replacements = @(
    @("2025-04-03 Thursday", "2025-04-04 Friday"),
    @("782×8=", "586×6="),
    @("360×2=", "257×8="),
    @("475×4=", "336×3="),
    @("720×7=", "808×7="),
    @("228×2=", "736×9="),
    @("403×8=", "428×4="),
    @("543×5=", "954×9="),
    @("946×3=", "115×6="),
    @("991×7=", "216×2="),
    @("204×6=", "302×7="),
    @("799×8=", "458×5="),
    @("370×8=", "589×6="),
    @("453×3=", "289×8="),
    @("982×9=", "833×6="),
    @("857×7=", "374×9="),
    @("438×4=", "265×2="),
    @("568×7=", "913×3="),
    @("726×8=", "421×3="),
    @("911×5=", "395×9="),
    @("735×6=", "771×5="),
    @("617×8=", "251×9="),
    @("236×9=", "558×8="),
    @("169×7=", "544×3="),
    @("162×2=", "166×7="),
    @("734×8=", "964×4=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
